# Update the build/version timestamp strings across the workbook.
$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value = "Version: mines - January 30 (built on " + $newStamp + ")"
$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Junde Coal Mine, China, M1156, version 'mines - January 30 (built on " + $newStamp + ")'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
$newBuildVersion = "mines - January 30 (built on " + $newStamp + ")"
for ($r = 2; $r -le 12; $r++) {
    $wsData.Range("S" + $r).Value = $newBuildVersion
}
